$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02,1.023304355680722,1.024755990762698,1.023621322831219,1.030961747353864,1.025953971324738,1.028485325793956,1.027583560050094,1.02645222249884,1.033771243027158,1.029945892925503),
    @(1.02,1.026377924265348,1.027570347497469,1.026314281422187,1.033695583306108,1.026035695987697,1.031186329626498,1.030199629062849,1.028946968458538,1.036308386270808,1.032650732493293),
    @(1.02,1.028351642016768,1.029377234708941,1.028042987718856,1.035448482401623,1.026084943210484,1.032919445564072,1.031878063888969,1.030547242842385,1.037933842631588,1.034386309654289),
    @(1.02,1.029177886260471,1.030133548300643,1.02876651900724,1.03618164214123,1.026104781800486,1.033644642379875,1.032580336563051,1.03121673533075,1.038613380593452,1.035112536332752),
    @(1.02,1.029316413682374,1.030260345793202,1.028887817072206,1.036304525020398,1.026108062210883,1.03376620911936,1.03269805790451,1.031328957226578,1.03872725738656,1.035234275710951),
    @(1.02,1.028362695973245,1.029387353456362,1.028052668099066,1.035458293570295,1.026085211686568,1.032929148935927,1.031887460700081,1.030556201339764,1.037942937485345,1.034396026806046),
    @(1.02,1.02434627039831,1.025710111947005,1.024534337488439,1.031889043111982,1.025982346125041,1.029401225646124,1.028470696753962,1.027298284022155,1.034632098024787,1.030863093460513),
    @(1.02,1.017148443357081,1.019117289724561,1.018224570409334,1.025472335588856,1.025773015367911,1.023068350306932,1.022335942857671,1.021446245542153,1.028669642804418,1.024521224711948),
    @(1.02,1.012262056575343,1.014639804573557,1.01393809484108,1.021103099696556,1.025614247915147,1.018762143590682,1.018163497759988,1.017464415314611,1.024602781357983,1.020208902688232),
    @(1.02,1.010123779910487,1.012680054577229,1.012061662181994,1.019188134417916,1.025540862763521,1.016876105583028,1.016335818013175,1.015719838417401,1.022818702131571,1.018320186290276),
    @(1.02,1.009326018064798,1.011948839280888,1.011361491515488,1.018473247778887,1.025512900300163,1.01617220441077,1.015653663059217,1.015068641552628,1.022152431359983,1.017615285497614),
    @(1.02,1.009497301845575,1.012105837807384,1.011511826392006,1.018626757497612,1.02551893031707,1.016323346905687,1.015800137720754,1.015208471443684,1.022295512765667,1.017766642632211),
    @(1.02,1.010057909108153,1.012619679687176,1.012003851443029,1.01912911546681,1.025538565780497,1.016817989816017,1.016279498359905,1.015666075903814,1.022763701771548,1.018261987992209),
    @(1.02,1.01040284829514,1.012935836555204,1.012306579241695,1.019438156342269,1.025550570337902,1.017122308914976,1.016574410844008,1.015947595766584,1.02305169005297,1.018566739259195),
    @(1.02,1.012403482388758,1.014769414222678,1.014062188277261,1.021229693760052,1.025619019926909,1.018886851868282,1.018284342526151,1.017579757157237,1.024720688489318,1.02033378806589),
    @(1.02,1.013652322690221,1.015913866935664,1.01515790009534,1.02234722182933,1.025660709898166,1.019987881277438,1.019351235565192,1.018598021824013,1.025761343181507,1.021436381063098),
    @(1.02,1.014378594722229,1.016579392344883,1.015795054145993,1.022996840591537,1.025684579972083,1.020628034775723,1.019971519974757,1.019189996325712,1.026366118365486,1.022077443652776),
    @(1.02,1.014625873286017,1.01680598138987,1.016011979032947,1.023217971406965,1.025692643447821,1.020845965153528,1.020182682540602,1.019391514986951,1.026571957726162,1.022295683516715),
    @(1.02,1.013518557880906,1.015791287260924,1.01504054378043,1.022227551751431,1.025656283243179,1.019869965030166,1.019236977450287,1.018488975431174,1.025649921329332,1.021318297361228),
    @(1.02,1.009892922296314,1.012468457554481,1.011859051104763,1.018981283457426,1.025532803117369,1.016672423099287,1.016138429948013,1.015531411618825,1.022625931653298,1.018116214554045),
    @(1.02,1.007592957184883,1.010360237946963,1.009840258897975,1.016919432199577,1.025451089664297,1.014642597391404,1.014171249072994,1.013653396066155,1.020703832652238,1.016083506260832),
    @(1.02,1.008814192145682,1.01147969167195,1.010912250239421,1.018014470355108,1.025494796394572,1.015720528240352,1.015215931313222,1.014650758750639,1.021724783882005,1.017162967895226),
    @(1.02,1.013579007092239,1.015846681956855,1.01509357808482,1.02228163235042,1.025658284838829,1.019923252659521,1.019288611997271,1.018538254858121,1.025700274924807,1.02137166066513),
    @(1.02,1.019024216851104,1.020835722515038,1.01986944020637,1.02714688100957,1.025830493025523,1.024719944956365,1.023936035037029,1.022972891004643,1.03022685523408,1.026175164815224)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$row"
        $ws.Range($addr).Value = $data[$i][$j]
    }
}
